# Fix typos in the "Tasa de serorreversion" / "Tasa de serorreversión Rhat"
# header labels (remove the extra "r": serorreversion -> seroreversion) and
# rename axis names, as described in the commit message
# "cambio de de debut sexual y nombres de ejes en español".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F8").Value = "Tasa de seroreversion"
$ws.Range("G8").Value = "Tasa de seroreversión Rhat"
